# Auto-generated edit script applying numeric corrections to the
# "Leve profit" data tables across all 8 sheets (ALC, ARM, BSM, CRP,
# CUL, GSM, LTW, WVR), per the scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 43407.6
$ws.Range("I21").Value = 43407.6
$ws.Range("K21").Value = 43407.6
$ws.Range("M21").Value = -42939.6

$ws.Range("H23").Value = 43407.6
$ws.Range("I23").Value = 43407.6
$ws.Range("K23").Value = 43407.6
$ws.Range("M23").Value = -43173.6

$ws.Range("H76").Value = 3351929.5
$ws.Range("I76").Value = 4690084.5
$ws.Range("J76").Value = 6542.5
$ws.Range("K76").Value = 4690084.5
$ws.Range("L76").Value = 6542.5
$ws.Range("M76").Value = -4689769.5
$ws.Range("N76").Value = -7172.5

$ws.Range("H79").Value = 3351929.5
$ws.Range("I79").Value = 4690084.5
$ws.Range("J79").Value = 6542.5
$ws.Range("K79").Value = 4690084.5
$ws.Range("L79").Value = 6542.5
$ws.Range("M79").Value = -4688992.5
$ws.Range("N79").Value = -8726.5

$ws.Range("H92").Value = 3078439.8
$ws.Range("I92").Value = 4104254
$ws.Range("K92").Value = 4104254
$ws.Range("M92").Value = -4103006

$ws.Range("H100").Value = 4250
$ws.Range("I100").Value = 5000
$ws.Range("K100").Value = 5000
$ws.Range("M100").Value = -4459

$ws.Range("H129").Value = 873.8163500000001
$ws.Range("J129").Value = 873.11365
$ws.Range("L129").Value = 2619.34095
$ws.Range("N129").Value = -12619.34095

$ws.Range("H135").Value = 534.9
$ws.Range("I135").Value = 483.22223
$ws.Range("K135").Value = 4349.00007
$ws.Range("M135").Value = -1814.00007

$ws.Range("H137").Value = 1827.0952
$ws.Range("I137").Value = 1617.3636
$ws.Range("K137").Value = 4852.0908
$ws.Range("M137").Value = -2302.0908

$ws.Range("H138").Value = 2953
$ws.Range("I138").Value = 4819.4
$ws.Range("J138").Value = 2369.75
$ws.Range("K138").Value = 14458.2
$ws.Range("L138").Value = 7109.25
$ws.Range("M138").Value = -9318.199999999999
$ws.Range("N138").Value = -17389.25

$ws.Range("H141").Value = 802366.1
$ws.Range("I141").Value = 1001852.1
$ws.Range("K141").Value = 3005556.3
$ws.Range("M141").Value = -3000376.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1543.037
$ws.Range("I61").Value = 1093.381
$ws.Range("K61").Value = 1093.381
$ws.Range("M61").Value = -881.3810000000001

$ws.Range("H74").Value = 1634.619
$ws.Range("J74").Value = 1727.6364
$ws.Range("L74").Value = 1727.6364
$ws.Range("N74").Value = -3475.6364

$ws.Range("H77").Value = 1634.619
$ws.Range("J77").Value = 1727.6364
$ws.Range("L77").Value = 8638.182000000001
$ws.Range("N77").Value = -17374.182

$ws.Range("H88").Value = 4466.1113
$ws.Range("I88").Value = 3265.3333
$ws.Range("J88").Value = 5066.5
$ws.Range("K88").Value = 3265.3333
$ws.Range("L88").Value = 5066.5
$ws.Range("M88").Value = -2859.3333
$ws.Range("N88").Value = -5878.5

$ws.Range("H91").Value = 4466.1113
$ws.Range("I91").Value = 3265.3333
$ws.Range("J91").Value = 5066.5
$ws.Range("K91").Value = 3265.3333
$ws.Range("L91").Value = 5066.5
$ws.Range("M91").Value = -1861.3333
$ws.Range("N91").Value = -7874.5

$ws.Range("H102").Value = 2393.1333
$ws.Range("I102").Value = 2036.1818
$ws.Range("K102").Value = 2036.1818
$ws.Range("M102").Value = -414.1818000000001

$ws.Range("H132").Value = 1419.7451
$ws.Range("I132").Value = 1044.0769
$ws.Range("K132").Value = 3132.2307
$ws.Range("M132").Value = -602.2307000000001

$ws.Range("H136").Value = 1543.037
$ws.Range("I136").Value = 1093.381
$ws.Range("K136").Value = 3280.143
$ws.Range("M136").Value = -730.143

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1287.7693
$ws.Range("I20").Value = 1313.9131
$ws.Range("K20").Value = 1313.9131
$ws.Range("M20").Value = -1066.9131

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2660.3635
$ws.Range("I99").Value = 1616.6666
$ws.Range("J99").Value = 3051.75
$ws.Range("K99").Value = 1616.6666
$ws.Range("L99").Value = 3051.75
$ws.Range("M99").Value = -118.6666
$ws.Range("N99").Value = -6047.75

$ws.Range("H122").Value = 3782.818
$ws.Range("I122").Value = 3379.8
$ws.Range("J122").Value = 4118.6665
$ws.Range("K122").Value = 10139.4
$ws.Range("L122").Value = 12355.9995
$ws.Range("M122").Value = -7689.400000000001
$ws.Range("N122").Value = -17255.9995

$ws.Range("H126").Value = 2660.3635
$ws.Range("I126").Value = 1616.6666
$ws.Range("J126").Value = 3051.75
$ws.Range("K126").Value = 4849.9998
$ws.Range("L126").Value = 9155.25
$ws.Range("M126").Value = -2379.9998
$ws.Range("N126").Value = -14095.25

$ws.Range("H132").Value = 2004.8223
$ws.Range("I132").Value = 1375.5358
$ws.Range("J132").Value = 3041.2942
$ws.Range("K132").Value = 4126.607400000001
$ws.Range("L132").Value = 9123.882599999999
$ws.Range("M132").Value = -1596.607400000001
$ws.Range("N132").Value = -14183.8826

$ws.Range("H134").Value = 1972.36
$ws.Range("I134").Value = 1686.1904
$ws.Range("K134").Value = 5058.5712
$ws.Range("M134").Value = -2523.5712

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1107.625
$ws.Range("I122").Value = 280
$ws.Range("K122").Value = 2520
$ws.Range("M122").Value = -70

$ws.Range("H131").Value = 12221.629
$ws.Range("J131").Value = 13513.619
$ws.Range("L131").Value = 40540.857
$ws.Range("N131").Value = -50620.857

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1097.5
$ws.Range("J80").Value = 1200
$ws.Range("L80").Value = 1200
$ws.Range("N80").Value = -3196

$ws.Range("H83").Value = 1097.5
$ws.Range("J83").Value = 1200
$ws.Range("L83").Value = 6000
$ws.Range("N83").Value = -15984

$ws.Range("H102").Value = 2962.5
$ws.Range("I102").Value = 2958.3333
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 2958.3333
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -1336.3333
$ws.Range("N102").Value = -6244

$ws.Range("H122").Value = 1882.9445
$ws.Range("I122").Value = 1300.75
$ws.Range("K122").Value = 3902.25
$ws.Range("M122").Value = -1452.25

$ws.Range("H126").Value = 2771.7407
$ws.Range("I126").Value = 2796.88
$ws.Range("J126").Value = 2457.5
$ws.Range("K126").Value = 8390.639999999999
$ws.Range("L126").Value = 7372.5
$ws.Range("M126").Value = -5920.639999999999
$ws.Range("N126").Value = -12312.5

$ws.Range("H132").Value = 1823.3455
$ws.Range("I132").Value = 1453.8536
$ws.Range("J132").Value = 2905.4285
$ws.Range("K132").Value = 4361.560799999999
$ws.Range("L132").Value = 8716.2855
$ws.Range("M132").Value = -1831.560799999999
$ws.Range("N132").Value = -13776.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2930.5

$ws.Range("H122").Value = 7374
$ws.Range("I122").Value = 5400.4
$ws.Range("J122").Value = 10663.333
$ws.Range("K122").Value = 16201.2
$ws.Range("L122").Value = 31989.999
$ws.Range("M122").Value = -13751.2
$ws.Range("N122").Value = -36889.999

$ws.Range("H132").Value = 1914.4103
$ws.Range("I132").Value = 1691.25
$ws.Range("K132").Value = 5073.75
$ws.Range("M132").Value = -2543.75

$ws.Range("H136").Value = 3427
$ws.Range("I136").Value = 2830.3333
$ws.Range("J136").Value = 3964
$ws.Range("K136").Value = 8490.999899999999
$ws.Range("L136").Value = 11892
$ws.Range("M136").Value = -5940.999899999999
$ws.Range("N136").Value = -16992

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2024.2858
$ws.Range("I81").Value = 2061.6667
$ws.Range("J81").Value = 1800
$ws.Range("K81").Value = 4123.3334
$ws.Range("L81").Value = 3600
$ws.Range("M81").Value = -3062.3334
$ws.Range("N81").Value = -5722

$ws.Range("H84").Value = 2024.2858
$ws.Range("I84").Value = 2061.6667
$ws.Range("J84").Value = 1800
$ws.Range("K84").Value = 20616.667
$ws.Range("L84").Value = 18000
$ws.Range("M84").Value = -15312.667
$ws.Range("N84").Value = -28608

$ws.Range("H113").Value = 1950
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1950
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 5850
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -10190

$ws.Range("H122").Value = 71888.73
$ws.Range("I122").Value = 87597.44500000001
$ws.Range("K122").Value = 262792.335
$ws.Range("M122").Value = -260342.335

$ws.Range("H132").Value = 1426.1875
$ws.Range("I132").Value = 1165
$ws.Range("K132").Value = 3495
$ws.Range("M132").Value = -965

$ws.Range("H136").Value = 2338.5454
$ws.Range("I136").Value = 2103.2666
$ws.Range("K136").Value = 6309.7998
$ws.Range("M136").Value = -3759.7998
